$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new blank column before column N (14th column: "Late"). This
# shifts the existing "Late" (col N) and "Outstanding" (col P) columns,
# together with all their data rows, one position to the right (N->O,
# P->Q), leaving a new blank column behind at N / O1, matching the
# author's manual column insert (RBI / Variable Instalments layout change).
$ws.Range("N1").EntireColumn.Insert()

# Restore the workbook's last active selection/cursor position.
$ws.Range("R10").Select()
